$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Y: header + empty cells for existing rows 2-7
$ws.Range("Y1").Value = "IDF Overlap"
$ws.Range("X1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("Y2").Value = ""
$ws.Range("Y3").Value = ""
$ws.Range("Y4").Value = ""
$ws.Range("Y5").Value = ""
$ws.Range("Y6").Value = ""
$ws.Range("Y7").Value = ""

# New row 8: Gemma-7B-Instruct
$ws.Range("A8").Value = "Gemma-7B-Instruct"
$ws.Range("B8").Value = "0.81 ± 0.39"
$ws.Range("C8").Value = "-0.03 ± 0.59"
$ws.Range("D8").Value = "0.3 ± 0.66"
$ws.Range("E8").Value = "0.01 ± 0.01"
$ws.Range("F8").Value = "0.12 ± 0.08"
$ws.Range("G8").Value = "0.02 ± 0.03"
$ws.Range("H8").Value = "0.11 ± 0.07"
$ws.Range("I8").Value = "0.13 ± 0.08"
$ws.Range("J8").Value = "0.83 ± 0.13"
$ws.Range("K8").Value = "0.84 ± 0.13"
$ws.Range("L8").Value = "0.83 ± 0.13"
$ws.Range("M8").Value = "0.86 ± 0.15"
$ws.Range("N8").Value = "0.97 ± 0.15"
$ws.Range("O8").Value = ""
$ws.Range("P8").Value = "0.47 ± 0.16"
$ws.Range("Q8").Value = "4.59 ± 0.65"
$ws.Range("R8").Value = "0.022 ± 0.00"
$ws.Range("S8").Value = "0.93 ± 0.14"
$ws.Range("T8").Value = "0.97 ± 0.15"
$ws.Range("U8").Value = "3.0 ± 1.15"
$ws.Range("V8").Value = "0.4 ± 0.42"
$ws.Range("W8").Value = "0.93 ± 0.14"
$ws.Range("X8").Value = "1.25 ± 0.3"
$ws.Range("Y8").Value = "0.07 ± 0.1"
